$wb = $excel.ActiveWorkbook

# This script applies updated market-price / profit figures pulled by the
# scheduled data-refresh runner to the per-sheet Leve profit tables.
# Each worksheet (ALC, ARM, BSM, CRP, CUL, GSM, LTW, WVR) holds the same
# column layout: H=currentAveragePrice, I=currentAveragePriceNQ,
# J=currentAveragePriceHQ, K=LevePriceNQ, L=LevePriceHQ, M=LeveProfitNQ,
# N=LeveProfitHQ.

# ----- Sheet: ALC -----
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H92").Value = 872.5
$ws.Range("I92").Value = 816.9231
$ws.Range("K92").Value = 816.9231
$ws.Range("M92").Value = 431.0769
$ws.Range("H129").Value = 832.3283699999999
$ws.Range("I129").Value = 599.2
$ws.Range("J129").Value = 851.129
$ws.Range("K129").Value = 1797.6
$ws.Range("L129").Value = 2553.387
$ws.Range("M129").Value = 3202.4
$ws.Range("N129").Value = -12553.387
$ws.Range("H137").Value = 92845.82000000001
$ws.Range("I137").Value = 0
$ws.Range("K137").Value = 0
$ws.Range("M137").ClearContents()
$ws.Range("H138").Value = 1591.9642
$ws.Range("I138").Value = 589.35297
$ws.Range("J138").Value = 3141.4546
$ws.Range("K138").Value = 1768.05891
$ws.Range("L138").Value = 9424.363799999999
$ws.Range("M138").Value = 3371.94109
$ws.Range("N138").Value = -19704.3638

# ----- Sheet: ARM -----
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H32").Value = 28360.512
$ws.Range("I32").Value = 29809.918
$ws.Range("K32").Value = 29809.918
$ws.Range("M32").Value = -29522.918
$ws.Range("H45").Value = 3374.8372
$ws.Range("I45").Value = 2919.3333
$ws.Range("K45").Value = 2919.3333
$ws.Range("M45").Value = -2542.3333
$ws.Range("H74").Value = 2399.4194
$ws.Range("I74").Value = 2521.9546
$ws.Range("J74").Value = 2099.889
$ws.Range("K74").Value = 2521.9546
$ws.Range("L74").Value = 2099.889
$ws.Range("M74").Value = -1647.9546
$ws.Range("N74").Value = -3847.889
$ws.Range("H77").Value = 2399.4194
$ws.Range("I77").Value = 2521.9546
$ws.Range("J77").Value = 2099.889
$ws.Range("K77").Value = 12609.773
$ws.Range("L77").Value = 10499.445
$ws.Range("M77").Value = -8241.773000000001
$ws.Range("N77").Value = -19235.445
$ws.Range("H96").Value = 18512.572
$ws.Range("J96").Value = 18512.572
$ws.Range("L96").Value = 18512.572
$ws.Range("N96").Value = -24004.572
$ws.Range("H97").Value = 735
$ws.Range("I97").Value = 726.3158
$ws.Range("K97").Value = 726.3158
$ws.Range("M97").Value = -230.3158
$ws.Range("H101").Value = 28457.572
$ws.Range("J101").Value = 28457.572
$ws.Range("L101").Value = 28457.572
$ws.Range("N101").Value = -34947.572
$ws.Range("H102").Value = 1481
$ws.Range("I102").Value = 1332.7142
$ws.Range("J102").Value = 2000
$ws.Range("K102").Value = 1332.7142
$ws.Range("L102").Value = 2000
$ws.Range("M102").Value = 289.2858000000001
$ws.Range("N102").Value = -5244
$ws.Range("H132").Value = 17763.562
$ws.Range("I132").Value = 2127.0435
$ws.Range("J132").Value = 57723.555
$ws.Range("K132").Value = 6381.130500000001
$ws.Range("L132").Value = 173170.665
$ws.Range("M132").Value = -3851.130500000001
$ws.Range("N132").Value = -178230.665
$ws.Range("H139").Value = 47903.168
$ws.Range("J139").Value = 47903.168
$ws.Range("L139").Value = 47903.168
$ws.Range("N139").Value = -58183.168

# ----- Sheet: BSM -----
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H81").Value = 20084.857
$ws.Range("J81").Value = 20084.857
$ws.Range("L81").Value = 20084.857
$ws.Range("N81").Value = -22206.857
$ws.Range("H84").Value = 20084.857
$ws.Range("J84").Value = 20084.857
$ws.Range("L84").Value = 60254.571
$ws.Range("N84").Value = -70862.571
$ws.Range("H94").Value = 1201.5
$ws.Range("I94").Value = 875.0526
$ws.Range("K94").Value = 875.0526
$ws.Range("M94").Value = -424.0526
$ws.Range("H99").Value = 1638.1666
$ws.Range("I99").Value = 1351.3889
$ws.Range("K99").Value = 1351.3889
$ws.Range("M99").Value = 146.6111000000001
$ws.Range("H134").Value = 36342.934
$ws.Range("I134").Value = 45123.918
$ws.Range("K134").Value = 135371.754
$ws.Range("M134").Value = -132836.754
$ws.Range("H137").Value = 36256.668
$ws.Range("J137").Value = 36256.668
$ws.Range("L137").Value = 36256.668
$ws.Range("N137").Value = -46456.668

# ----- Sheet: CRP -----
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H31").Value = 14820.125
$ws.Range("I31").Value = 17921.277
$ws.Range("J31").Value = 5516.6665
$ws.Range("K31").Value = 17921.277
$ws.Range("L31").Value = 5516.6665
$ws.Range("M31").Value = -17626.277
$ws.Range("N31").Value = -6106.6665
$ws.Range("H34").Value = 14820.125
$ws.Range("I34").Value = 17921.277
$ws.Range("J34").Value = 5516.6665
$ws.Range("K34").Value = 17921.277
$ws.Range("L34").Value = 5516.6665
$ws.Range("M34").Value = -17719.277
$ws.Range("N34").Value = -5920.6665
$ws.Range("H134").Value = 1236.1578
$ws.Range("I134").Value = 971.2
$ws.Range("K134").Value = 2913.6
$ws.Range("M134").Value = -378.6000000000004

# ----- Sheet: CUL -----
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H5").Value = 2051.2222
$ws.Range("I5").Value = 1779.4286
$ws.Range("K5").Value = 5338.2858
$ws.Range("M5").Value = -5226.2858
$ws.Range("H122").Value = 467.75
$ws.Range("I122").Value = 249.14285
$ws.Range("K122").Value = 2242.28565
$ws.Range("M122").Value = 207.7143499999997
$ws.Range("H131").Value = 770.567
$ws.Range("J131").Value = 793.5955
$ws.Range("L131").Value = 2380.7865
$ws.Range("N131").Value = -12460.7865
$ws.Range("H135").Value = 2051.2222
$ws.Range("I135").Value = 1779.4286
$ws.Range("K135").Value = 16014.8574
$ws.Range("M135").Value = -13479.8574

# ----- Sheet: GSM -----
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H97").Value = 1793
$ws.Range("I97").Value = 986.8570999999999
$ws.Range("J97").Value = 3674
$ws.Range("K97").Value = 986.8570999999999
$ws.Range("L97").Value = 3674
$ws.Range("M97").Value = -490.8570999999999
$ws.Range("N97").Value = -4666
$ws.Range("H132").Value = 119792.234
$ws.Range("J132").Value = 86583
$ws.Range("L132").Value = 259749
$ws.Range("N132").Value = -264809

# ----- Sheet: LTW -----
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H22").Value = 4628.857
$ws.Range("I22").Value = 3633.6667
$ws.Range("J22").Value = 5375.25
$ws.Range("K22").Value = 3633.6667
$ws.Range("L22").Value = 5375.25
$ws.Range("M22").Value = -3338.6667
$ws.Range("N22").Value = -5965.25
$ws.Range("H27").Value = 4628.857
$ws.Range("I27").Value = 3633.6667
$ws.Range("J27").Value = 5375.25
$ws.Range("K27").Value = 3633.6667
$ws.Range("L27").Value = 5375.25
$ws.Range("M27").Value = -3526.6667
$ws.Range("N27").Value = -5589.25
$ws.Range("H46").Value = 1100.25
$ws.Range("I46").Value = 1100.5
$ws.Range("J46").Value = 1100
$ws.Range("K46").Value = 1100.5
$ws.Range("L46").Value = 1100
$ws.Range("M46").Value = -912.5
$ws.Range("N46").Value = -1476
$ws.Range("H55").Value = 392.4
$ws.Range("I55").Value = 280
$ws.Range("J55").Value = 420.5
$ws.Range("K55").Value = 280
$ws.Range("L55").Value = 420.5
$ws.Range("M55").Value = -107
$ws.Range("N55").Value = -766.5
$ws.Range("H61").Value = 4083.2856
$ws.Range("J61").Value = 5163
$ws.Range("L61").Value = 5163
$ws.Range("N61").Value = -5567
$ws.Range("H100").Value = 1824.4
$ws.Range("I100").Value = 1616.9166
$ws.Range("K100").Value = 1616.9166
$ws.Range("M100").Value = -1075.9166
$ws.Range("H113").Value = 4083.2856
$ws.Range("J113").Value = 5163
$ws.Range("L113").Value = 5163
$ws.Range("N113").Value = -9503
$ws.Range("H122").Value = 1786257.2
$ws.Range("I122").Value = 2453953
$ws.Range("J122").Value = 5735
$ws.Range("K122").Value = 7361859
$ws.Range("L122").Value = 17205
$ws.Range("M122").Value = -7359409
$ws.Range("N122").Value = -22105

# ----- Sheet: WVR -----
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H54").Value = 19666.666
$ws.Range("J54").Value = 19666.666
$ws.Range("L54").Value = 19666.666
$ws.Range("N54").Value = -20706.666
$ws.Range("H132").Value = 1826.909
$ws.Range("I132").Value = 933.1667
$ws.Range("K132").Value = 2799.5001
$ws.Range("M132").Value = -269.5001000000002
$ws.Range("H136").Value = 25001716
$ws.Range("I136").Value = 43479892
$ws.Range("J136").Value = 1829.4706
$ws.Range("K136").Value = 130439676
$ws.Range("L136").Value = 5488.4118
$ws.Range("M136").Value = -130437126
$ws.Range("N136").Value = -10588.4118
